$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add a new BOM row 29 for the new opto-relay parts (CPC1018N) ---
# Seed formatting by copying the last existing data row (28) down into the
# new row, then overwrite the values.
$ws.Range("A28:D28").Copy()
$ws.Range("A29:D29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(29).RowHeight = 20

$ws.Cells.Item(29, 3).Value = "SOP-4-2.54mm"
$ws.Cells.Item(29, 2).Value = "U7, U8"

# --- Update existing BOM rows whose designator lists changed ---

# Row 22 (1K / R19, R20, R21) -> add R24, R25
$ws.Range("B22").Value = "R19, R20, R21, R24, R25"

# Row 23 (10K / R22-29) -> now explicit designators for the remaining parts
$ws.Range("B23").Value = "R22, R23, R26, R27"

# Row 5 (10nF / C18-21) -> split into C18, C19
$ws.Range("B5").Value = "C18, C19"
$ws.Range("B5").Font.Name = "Arial"
$ws.Range("B5").Font.Size = 11
$ws.Range("B5").Font.Color = 0

$ws.Cells.Item(29, 1).Value = "CPC1018N"
$ws.Cells.Item(29, 4).Value = "C1558973"

# New row's designator cell has no border/fill, unlike the copied row above it
$ws.Range("B29").Borders.LineStyle = -4142
$ws.Range("B29").Interior.Pattern = -4142

# Match the final selection left behind in the saved workbook
$ws.Range("D29").Select()
